# Apply the edit: rows of MAG data were re-derived/re-sorted and the
# table was shortened from 21 data rows (A1:F21) down to 15 data rows
# (A1:F15). Rows 2-3 (GUT11063, GUT11107) are unchanged; rows 4-15 get
# new content, and the former rows 16-21 are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 4 through 15 (columns A-F).
# NOTE: values are written in plain decimal form (no exponent) because
# the interpreter's numeric literal grammar does not accept scientific
# notation (e.g. "1e-08"); the decimal expansions below are exact.
$data = @(
    @("even_MAG-GUT23816.fa", 0.9999999830990612, 0.00000001690093876783204, 0.9999999830990612, "f__Anaerovoracaceae", "f__Anaerovoracaceae"),
    @("even_MAG-GUT25429.fa", 0.0000003628757221374457, 0.9999996371242779, 0.9999996371242779, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT30539.fa", 0.00000009588940319549266, 0.9999999041105968, 0.9999999041105968, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT32544.fa", 0.0000005946883073759679, 0.9999994053116926, 0.9999994053116926, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT42774.fa", 0.999999959521361, 0.00000004047863894791822, 0.999999959521361, "f__Anaerovoracaceae", "f__Anaerovoracaceae"),
    @("even_MAG-GUT44190.fa", 0.00000007769526022549655, 0.9999999223047398, 0.9999999223047398, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT44598.fa", 0.000003228625496065618, 0.9999967713745039, 0.9999967713745039, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT56457.fa", 0.0000004533934661044015, 0.9999995466065339, 0.9999995466065339, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT58695.fa", 0.0000002061700127287125, 0.9999997938299873, 0.9999997938299873, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT66161.fa", 0.0000025182906178145, 0.9999974817093822, 0.9999974817093822, "f__Peptostreptococcaceae", "f__Peptostreptococcaceae"),
    @("even_MAG-GUT85090.fa", 0.9999934872148809, 0.000006512785119105472, 0.9999934872148809, "f__Anaerovoracaceae", "f__Anaerovoracaceae"),
    @("even_MAG-GUT85141.fa", 0.9999999916410939, 0.000000008358906107604512, 0.9999999916410939, "f__Anaerovoracaceae", "f__Anaerovoracaceae")
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}

# Remove the now-unused trailing rows (16-21), which are no longer part
# of the (now 15-row) table.
$ws.Range("A16:F21").Delete()
